$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a cell with default (unstyled) formatting as the format source so that
# forcing text values (via a leading apostrophe) does not leave a stray
# quote-prefixed / text-numfmt style on the edited cells.
$fmtSrc = $ws.Range("B2")
$fmtSrc.Copy() | Out-Null

$c = $ws.Range("D2")
$c.Value = "'47.379.30"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D3")
$c.Value = "'2.489.55"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E3")
$c.Value = "'  +0.01%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E4")
$c.Value = "'  -0.01%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D5")
$c.Value = "'321.52"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E5")
$c.Value = "'  -0.33%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D6")
$c.Value = "'109.04"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E6")
$c.Value = "'  +3.41%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D7")
$c.Value = "'0.522"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E7")
$c.Value = "'  -0.36%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D8")
$c.Value = "'0.999"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E8")
$c.Value = "'  -0.04%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E9")
$c.Value = "'  -0.15%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D10")
$c.Value = "'39.39"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E10")
$c.Value = "'  +4.01%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E11")
$c.Value = "'  -0.40%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E12")
$c.Value = "'  +1.01%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D13")
$c.Value = "'18.60"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E13")
$c.Value = "'  +1.77%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E14")
$c.Value = "'  +0.58%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D15")
$c.Value = "'2.877.98"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E15")
$c.Value = "'  -0.13%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D16")
$c.Value = "'2.488.33"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E16")
$c.Value = "'  -0.01%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D17")
$c.Value = "'0.847"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E17")
$c.Value = "'  +0.56%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D18")
$c.Value = "'47.280.42"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E18")
$c.Value = "'  -0.13%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D19")
$c.Value = "'13.45"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E19")
$c.Value = "'  +6.01%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E20")
$c.Value = "'  +1.19%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D21")
$c.Value = "'0.0₃0941"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D22")
$c.Value = "'2.75"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E22")
$c.Value = "'  +15.75%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D23")
$c.Value = "'70.65"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E23")
$c.Value = "'  +0.00%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D24")
$c.Value = "'246.88"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E24")
$c.Value = "'  -1.42%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E25")
$c.Value = "'  -0.59%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E26")
$c.Value = "'  -0.02%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D27")
$c.Value = "'25.75"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E27")
$c.Value = "'  -1.56%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E28")
$c.Value = "'  +0.69%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D29")
$c.Value = "'9.98"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E29")
$c.Value = "'  -0.27%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E30")
$c.Value = "'  +4.18%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D31")
$c.Value = "'34.68"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E31")
$c.Value = "'  -0.85%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D32")
$c.Value = "'49.89"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E32")
$c.Value = "'  +0.87%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D33")
$c.Value = "'20.43"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E33")
$c.Value = "'  +2.72%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D34")
$c.Value = "'5.31"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E34")
$c.Value = "'  -0.54%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D35")
$c.Value = "'0.0787"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E35")
$c.Value = "'  +0.62%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E36")
$c.Value = "'  +0.09%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E37")
$c.Value = "'  +2.40%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E38")
$c.Value = "'  +0.91%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E39")
$c.Value = "'  -2.06%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E40")
$c.Value = "'  +0.31%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D41")
$c.Value = "'22.43"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E41")
$c.Value = "'  +7.12%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E42")
$c.Value = "'  -1.97%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D43")
$c.Value = "'119.10"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E43")
$c.Value = "'  -2.11%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E44")
$c.Value = "'  -0.33%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D45")
$c.Value = "'1.992.96"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E45")
$c.Value = "'  +1.70%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E46")
$c.Value = "'  +1.97%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D47")
$c.Value = "'2.04"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E47")
$c.Value = "'  -2.72%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("B48")
$c.Value = "'FraxShare"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("C48")
$c.Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D48")
$c.Value = "'9.09"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E48")
$c.Value = "'  -1.36%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("B49")
$c.Value = "'Stacks"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("C49")
$c.Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D49")
$c.Value = "'1.78"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E49")
$c.Value = "'  -0.79%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E50")
$c.Value = "'  -2.38%  "
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("D51")
$c.Value = "'56.80"
$c.PasteSpecial(-4122) | Out-Null

$c = $ws.Range("E51")
$c.Value = "'  +3.60%  "
$c.PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
